$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vocabulary rows (B66:C89), written in the same cell-by-cell order the
# original author used (this matters for the shared-string table order).
$ws.Cells.Item(66, 2).Value = "gotta (I have got to)"
$ws.Cells.Item(66, 3).Value = "мне нужно, я вынужден "
$ws.Cells.Item(67, 2).Value = "to part ways"
$ws.Cells.Item(67, 3).Value = "Закончиться, разойтись, расстаться"
$ws.Cells.Item(68, 2).Value = "sick and tired"
$ws.Cells.Item(69, 2).Value = "Mockingjay"
$ws.Cells.Item(69, 3).Value = "сойка-пересмешника"
$ws.Cells.Item(68, 3).Value = "Откровенно достало"
$ws.Cells.Item(71, 3).Value = "застегнуть (ремень)"
$ws.Cells.Item(70, 3).Value = "Застегнуть молнию"
$ws.Cells.Item(72, 2).Value = "to twist"
$ws.Cells.Item(72, 3).Value = "сделать человеку плохо"
$ws.Cells.Item(71, 2).Value = "to fasten"
$ws.Cells.Item(70, 2).Value = "to zip"
$ws.Cells.Item(73, 3).Value = "Притормозить (про человека)"
$ws.Cells.Item(74, 2).Value = "to let up"
$ws.Cells.Item(74, 3).Value = "Ослабевать, сходить на нет"
$ws.Cells.Item(73, 2).Value = "to slow up (down)"
$ws.Cells.Item(75, 2).Value = "sacrilegious"
$ws.Cells.Item(75, 3).Value = "кощунственный"
$ws.Cells.Item(76, 2).Value = "take forever"
$ws.Cells.Item(76, 3).Value = "Очень долго (время)"
$ws.Cells.Item(77, 2).Value = "run-in"
$ws.Cells.Item(77, 3).Value = "схватка"
$ws.Cells.Item(78, 2).Value = "to detract"
$ws.Cells.Item(78, 3).Value = "Уменьшать "
$ws.Cells.Item(79, 2).Value = "detractors"
$ws.Cells.Item(79, 3).Value = "недоброжелатели"
$ws.Cells.Item(80, 2).Value = "pitchfork"
$ws.Cells.Item(80, 3).Value = "вилы"
$ws.Cells.Item(81, 2).Value = "corkscrew"
$ws.Cells.Item(81, 3).Value = "штопор"
$ws.Cells.Item(82, 2).Value = "revenge"
$ws.Cells.Item(82, 3).Value = "расплата, месть"
$ws.Cells.Item(83, 2).Value = "leeches"
$ws.Cells.Item(83, 3).Value = "пиявки"
$ws.Cells.Item(84, 2).Value = "flip the bird"
$ws.Cells.Item(84, 3).Value = "Поднять средний палец"
$ws.Cells.Item(85, 2).Value = "Lip-syncs"
$ws.Cells.Item(85, 3).Value = "Липсинг"
$ws.Cells.Item(86, 2).Value = "hardly"
$ws.Cells.Item(86, 3).Value = "Едва ли"
$ws.Cells.Item(87, 2).Value = "prolly"
$ws.Cells.Item(87, 3).Value = "возможно (prodadly)"
$ws.Cells.Item(88, 2).Value = "to fetch"
$ws.Cells.Item(89, 2).Value = "remote"
$ws.Cells.Item(89, 3).Value = "пульт"
$ws.Cells.Item(88, 3).Value = "пойди принеси (часто говорят собакам)"

# Give column A an explicit width (matches the new <col min="1" max="1".../>
# entry); leave columns B:C untouched so their existing merged <col> entry
# survives as-is.
$ws.Columns.Item(1).ColumnWidth = 18.8

# Selection / zoom match the saved view state in the edited workbook.
[void]$ws.Range("B15").Select()
$excel.ActiveWindow.Zoom = 115
